$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-deleted last row (old row 25), shifting content below it up by one
$ws.Rows(25).Delete()

# Set final cell values for rows 1-24 to match the target layout
$ws.Range("A1").Value = $null
$ws.Range("B1").Value = "Ementa atual:"
$ws.Range("C1").Value = "Ementa modificada (dados modificados em vermelho):"

$ws.Range("A2").Value = $null
$ws.Range("B2").Value = "LOQ4252"
$ws.Range("C2").Value = "LOQ4252"

$ws.Range("A3").Value = "Nome:"
$ws.Range("B3").Value = " Fundamentos de Fenômenos de Transporte"
$ws.Range("C3").Value = " Fundamentos de Fenômenos de Transporte"

$ws.Range("A4").Value = "Name:"
$ws.Range("B4").Value = "Transport Phenomena Fundamentals"
$ws.Range("C4").Value = "Transport Phenomena Fundamentals"

$ws.Range("A5").Value = "Créditos-aula:"
$ws.Range("B5").Value = "4"
$ws.Range("C5").Value = "4"

$ws.Range("A6").Value = "Créditos-trabalho"
$ws.Range("B6").Value = "0"
$ws.Range("C6").Value = "0"

$ws.Range("A7").Value = "Carga horária:"
$ws.Range("B7").Value = "60 h"
$ws.Range("C7").Value = "60 h"

$ws.Range("A8").Value = "Ativação:"
$ws.Range("B8").Value = "01/01/2021"
$ws.Range("C8").Value = "01/01/2021"

$ws.Range("A9").Value = "Semestre ideal:"
$ws.Range("B9").Value = "EP-5"
$ws.Range("C9").Value = "EP-5"

$ws.Range("A10").Value = "Objetivos:"
$ws.Range("B10").Value = "4808662 - Lucrécio Fábio dos Santos"
$ws.Range("C10").Value = "4808662 - Lucrécio Fábio dos Santos"
$ws.Rows(10).RowHeight = 60

$ws.Range("A11").Value = "Objectives:"
$ws.Range("B11").Value = $null
$ws.Range("C11").Value = $null
$ws.Rows(11).RowHeight = 60

$ws.Range("A12").Value = "Docentes responsáveis:"
$ws.Range("B12").Value = $null
$ws.Range("C12").Value = $null

$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Rows(13).RowHeight = 60

$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = $null
$ws.Range("C14").Value = $null
$ws.Rows(14).RowHeight = 60

$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2021"
$ws.Range("C15").Value = "01/01/2021"
$ws.Rows(15).RowHeight = 120

$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = $null
$ws.Range("C16").Value = $null
$ws.Rows(16).RowHeight = 120

$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17").Value = $null
$ws.Range("C17").Value = $null

$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "4808662 - Lucrécio Fábio dos Santos"
$ws.Range("C18").Value = "4808662 - Lucrécio Fábio dos Santos"
$ws.Rows(18).RowHeight = 60

$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Aulas expositivas teóricas, aulas práticas, aulas de exercícios, aulas de laboratório."
$ws.Range("C19").Value = "Aulas expositivas teóricas, aulas práticas, aulas de exercícios, aulas de laboratório."
$ws.Rows(19).RowHeight = 60

$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "Nota de duas provas (P1 e P2)Fórmula: M1 = (P1 + 2 x P2)/3.."
$ws.Range("C20").Value = "Nota de duas provas (P1 e P2)Fórmula: M1 = (P1 + 2 x P2)/3.."
$ws.Rows(20).RowHeight = 60

$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "Aplicação de uma prova envolvendo o assunto de todo semestre.NR (nota da recuperação) = (M1 + NR)/2."
$ws.Range("C21").Value = "Aplicação de uma prova envolvendo o assunto de todo semestre.NR (nota da recuperação) = (M1 + NR)/2."
$ws.Rows(21).RowHeight = 120

$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22").Value = $null
$ws.Range("C22").Value = $null

$ws.Range("A23").Value = $null
$ws.Range("B23").Value = "LOB1004 -  Cálculo II  (Requisito fraco)`n"
$ws.Range("C23").Value = "LOB1004 -  Cálculo II  (Requisito fraco)`n"
$ws.Rows(23).RowHeight = 30

$ws.Range("A24").Value = $null
$ws.Range("B24").Value = "LOB1019 -  Física II  (Requisito fraco)`n"
$ws.Range("C24").Value = "LOB1019 -  Física II  (Requisito fraco)`n"
$ws.Rows(24).RowHeight = 30
